# Adds rows 47-51 to the "Artfynd" sheet with new species-finding records,
# matching the source data export (Swedish species observation records).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47 ---
$ws.Range("A47").Value = 111974187
$ws.Range("B47").Value = 90710
$ws.Range("C47").Value = 'Ovaliderad'
$ws.Range("D47").Value = 'NT'
$ws.Range("E47").Value = 5449
$ws.Range("F47").Value = 'Svart taggsvamp'
$ws.Range("G47").Value = 'Phellodon niger'
$ws.Range("H47").Value = '(Fr.:Fr.) P.Karst.'
$ws.Range("I47").NumberFormat = "@"
$ws.Range("I47").Value = '1'
$ws.Range("J47").Value = 'mycel'
$ws.Range("P47").Value = 'Gustavbacke, norr om ån och väg, Jmt'
$ws.Range("Q47").Value = 439865.3631964622
$ws.Range("R47").Value = 6952242.088420792
$ws.Range("S47").Value = 10
$ws.Range("T47").Value = 'Jämtland'
$ws.Range("U47").Value = 'Härjedalen'
$ws.Range("V47").Value = 'Jämtland'
$ws.Range("W47").Value = 'Vemdalen'
$ws.Range("Y47").NumberFormat = "@"
$ws.Range("Y47").Value = '2023-09-05'
$ws.Range("Z47").Value = '00:00'
$ws.Range("AA47").NumberFormat = "@"
$ws.Range("AA47").Value = '2023-09-05'
$ws.Range("AB47").Value = '00:00'
$ws.Range("AD47").Value = $false
$ws.Range("AE47").Value = $false
$ws.Range("AG47").Value = $false
$ws.Range("AI47").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark, under tallåga'
# AT47: empty string cell in source data - nothing to set
$ws.Range("AW47").Value = 'Magnus Andersson'
$ws.Range("AX47").Value = 'Magnus Andersson'
$ws.Range("AY47").Value = 'SCA Skog Naturvärdesinventering'

# --- Row 48 ---
$ws.Range("A48").Value = 111974191
$ws.Range("B48").Value = 90652
$ws.Range("C48").Value = 'Ovaliderad'
$ws.Range("D48").Value = 'NT'
$ws.Range("E48").Value = 3100
$ws.Range("F48").Value = 'Talltaggsvamp'
$ws.Range("G48").Value = 'Bankera fuligineoalba'
$ws.Range("H48").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("I48").NumberFormat = "@"
$ws.Range("I48").Value = '1'
$ws.Range("J48").Value = 'mycel'
$ws.Range("P48").Value = 'Gustavbacke, norr om ån och väg, Jmt'
$ws.Range("Q48").Value = 439977.5118376439
$ws.Range("R48").Value = 6952213.872195411
$ws.Range("S48").Value = 10
$ws.Range("T48").Value = 'Jämtland'
$ws.Range("U48").Value = 'Härjedalen'
$ws.Range("V48").Value = 'Jämtland'
$ws.Range("W48").Value = 'Vemdalen'
$ws.Range("Y48").NumberFormat = "@"
$ws.Range("Y48").Value = '2023-09-05'
$ws.Range("Z48").Value = '00:00'
$ws.Range("AA48").NumberFormat = "@"
$ws.Range("AA48").Value = '2023-09-05'
$ws.Range("AB48").Value = '00:00'
$ws.Range("AD48").Value = $false
$ws.Range("AE48").Value = $false
$ws.Range("AG48").Value = $false
$ws.Range("AI48").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT48: empty string cell in source data - nothing to set
$ws.Range("AW48").Value = 'Magnus Andersson'
$ws.Range("AX48").Value = 'Magnus Andersson'
$ws.Range("AY48").Value = 'SCA Skog Naturvärdesinventering'

# --- Row 49 ---
$ws.Range("A49").Value = 111974186
$ws.Range("B49").Value = 90682
$ws.Range("C49").Value = 'Ovaliderad'
$ws.Range("D49").Value = 'NT'
$ws.Range("E49").Value = 2059
$ws.Range("F49").Value = 'Skrovlig taggsvamp'
$ws.Range("G49").Value = 'Hydnellum scabrosum'
$ws.Range("H49").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("I49").NumberFormat = "@"
$ws.Range("I49").Value = '1'
$ws.Range("J49").Value = 'mycel'
$ws.Range("P49").Value = 'Gustavbacke, norr om ån och väg, Jmt'
$ws.Range("Q49").Value = 439860.448822267
$ws.Range("R49").Value = 6952249.98427855
$ws.Range("S49").Value = 10
$ws.Range("T49").Value = 'Jämtland'
$ws.Range("U49").Value = 'Härjedalen'
$ws.Range("V49").Value = 'Jämtland'
$ws.Range("W49").Value = 'Vemdalen'
$ws.Range("Y49").NumberFormat = "@"
$ws.Range("Y49").Value = '2023-09-05'
$ws.Range("Z49").Value = '00:00'
$ws.Range("AA49").NumberFormat = "@"
$ws.Range("AA49").Value = '2023-09-05'
$ws.Range("AB49").Value = '00:00'
$ws.Range("AD49").Value = $false
$ws.Range("AE49").Value = $false
$ws.Range("AG49").Value = $false
$ws.Range("AI49").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT49: empty string cell in source data - nothing to set
$ws.Range("AW49").Value = 'Magnus Andersson'
$ws.Range("AX49").Value = 'Magnus Andersson'
$ws.Range("AY49").Value = 'SCA Skog Naturvärdesinventering'

# --- Row 50 ---
$ws.Range("A50").Value = 111974185
$ws.Range("B50").Value = 90660
$ws.Range("C50").Value = 'Ovaliderad'
$ws.Range("D50").Value = 'NT'
$ws.Range("E50").Value = 4362
$ws.Range("F50").Value = 'Blå taggsvamp'
$ws.Range("G50").Value = 'Hydnellum caeruleum'
$ws.Range("H50").Value = '(Hornem.) P.Karst.'
$ws.Range("I50").NumberFormat = "@"
$ws.Range("I50").Value = '1'
$ws.Range("J50").Value = 'mycel'
$ws.Range("P50").Value = 'Gustavbacke, norr om ån och väg, Jmt'
$ws.Range("Q50").Value = 439827.4842555065
$ws.Range("R50").Value = 6952232.676732311
$ws.Range("S50").Value = 10
$ws.Range("T50").Value = 'Jämtland'
$ws.Range("U50").Value = 'Härjedalen'
$ws.Range("V50").Value = 'Jämtland'
$ws.Range("W50").Value = 'Vemdalen'
$ws.Range("Y50").NumberFormat = "@"
$ws.Range("Y50").Value = '2023-09-05'
$ws.Range("Z50").Value = '00:00'
$ws.Range("AA50").NumberFormat = "@"
$ws.Range("AA50").Value = '2023-09-05'
$ws.Range("AB50").Value = '00:00'
$ws.Range("AD50").Value = $false
$ws.Range("AE50").Value = $false
$ws.Range("AG50").Value = $false
$ws.Range("AI50").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT50: empty string cell in source data - nothing to set
$ws.Range("AW50").Value = 'Magnus Andersson'
$ws.Range("AX50").Value = 'Magnus Andersson'
$ws.Range("AY50").Value = 'SCA Skog Naturvärdesinventering'

# --- Row 51 ---
$ws.Range("A51").Value = 111974188
$ws.Range("B51").Value = 90652
$ws.Range("C51").Value = 'Ovaliderad'
$ws.Range("D51").Value = 'NT'
$ws.Range("E51").Value = 3100
$ws.Range("F51").Value = 'Talltaggsvamp'
$ws.Range("G51").Value = 'Bankera fuligineoalba'
$ws.Range("H51").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("I51").NumberFormat = "@"
$ws.Range("I51").Value = '1'
$ws.Range("J51").Value = 'mycel'
$ws.Range("P51").Value = 'Gustavbacke, norr om ån och väg, Jmt'
$ws.Range("Q51").Value = 439869.6589509377
$ws.Range("R51").Value = 6952225.479112641
$ws.Range("S51").Value = 10
$ws.Range("T51").Value = 'Jämtland'
$ws.Range("U51").Value = 'Härjedalen'
$ws.Range("V51").Value = 'Jämtland'
$ws.Range("W51").Value = 'Vemdalen'
$ws.Range("Y51").NumberFormat = "@"
$ws.Range("Y51").Value = '2023-09-05'
$ws.Range("Z51").Value = '00:00'
$ws.Range("AA51").NumberFormat = "@"
$ws.Range("AA51").Value = '2023-09-05'
$ws.Range("AB51").Value = '00:00'
$ws.Range("AD51").Value = $false
$ws.Range("AE51").Value = $false
$ws.Range("AG51").Value = $false
$ws.Range("AI51").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT51: empty string cell in source data - nothing to set
$ws.Range("AW51").Value = 'Magnus Andersson'
$ws.Range("AX51").Value = 'Magnus Andersson'
$ws.Range("AY51").Value = 'SCA Skog Naturvärdesinventering'

Write-Output "Added rows 47-51 to Artfynd sheet"
